$wb = $excel.ActiveWorkbook

# --- locaties sheet: remove duplicate rows ---
$wsLoc = $wb.Worksheets.Item("locaties")
$wsLoc.Rows.Item(3).Delete()
$wsLoc.Rows.Item(4).Delete()

# --- metingen sheet: move selection ---
$wsMet = $wb.Worksheets.Item("metingen")
$wsMet.Range("D1").Select()

# --- locaties sheet becomes the active tab with a new selection ---
$wsLoc.Activate()
$wsLoc.Range("B13").Select()
